# Target sheet: "老師上課 + 監考時數" (the 2nd worksheet) holds the per-teacher
# exam-proctoring hours table (rows 2-167) with a summary row (168).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# 1) Wrap every "G168*C<row>" average-hours formula in ROUND(..., 0) so the
#    per-teacher allotted hours are whole numbers.
$rows = @(2,7,11,14,16,18,20,24,26,27,29,33,36,37,40,43,44,46,50,54,57,59,
          64,68,70,76,77,79,82,85,90,94,96,99,102,105,108,110,112,114,116,
          118,120,121,122,123,125,130,132,133,135,138,142,145,149,151,154,
          158,161,162,163,164,166,167)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Formula = "=ROUND(G168*C$r, 0)"
}

# 2) Add column-total array formulas on row 168 for C:F, and turn the
#    existing G168 average formula into an (equivalent) array formula too.
$ws.Range("C168").FormulaArray = '=SUM(C2:C167)'
$ws.Range("D168").FormulaArray = '=SUM(D2:D167)'
$ws.Range("E168").FormulaArray = '=SUM(E2:E167)'
$ws.Range("F168").FormulaArray = '=SUM(F2:F167)'
$ws.Range("G168").FormulaArray = '=ROUND((SUM($D$2:$D$167*$C$2:$C$167)-F161-F162-F163-F164+SUM($E$2:$E$167))/SUM($C$2:$C$167),0)'
